# "Add files via upload" — refresh the VLOOKUP lookup table on Sheet3 with a
# new data pull, then snapshot the recalculated results into a brand-new
# dated column ("13-nov") on Sheet1, right after the previous snapshot
# ("10-nov" in column CM).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------
# 1) Overwrite the raw lookup values in Sheet3!B20:B36 (key/value table
#    that Sheet3!C2:C18 and, via VLOOKUP, Sheet1!CB:CC depend on).
#    Only the rows below actually change value; the rest stay as-is.
# ---------------------------------------------------------------------
$newLookup = @{
    20 = 12.946338812456732   # 3D QUESO 92GX27
    23 = 8.5091935859377461   # DORITOS QUESO 70X40G
    24 = 7.8164448984773891   # DORITOS QUESO 77GX26
    26 = 7.0469660803807077   # LAYS CLASICAS 145GRX18
    27 = 8.4338045207842978   # LAYS CLASICAS 249GRX14
    28 = 7.0225947813515264   # LAYS CLASICAS 40GX68
    29 = 6.828685610610159    # LAYS CLASICAS 94GRX25
    30 = 0.62496000000192742  # LAYS ONDAS FH 30GX72
    31 = 3.1997734999998424   # LAYS ONDAS FH 70GX28
    32 = 16.911388755920086   # LAYS QSO Y CEBOLLA 34GX72
    33 = 14.901707307121656   # PEHUAMAR ACANALADA 520GX9
    34 = 8.0197564110179638   # PEHUAMAR MAICITOS 285GX10
    35 = 6.7996246164144107   # PEHUAMAR PAPA LISA 520GX9
    36 = 45.433600870996599   # QUAKER AVENA INSTANT FORTIF 18X280G
}

foreach ($row in $newLookup.Keys) {
    $ws3.Range("B$row").Value = $newLookup[$row]
}

# ---------------------------------------------------------------------
# 2) Add the new dated column on Sheet1 (column CN, right after the
#    existing last column CM = "10-nov"). Header is a new shared string
#    "13-nov"; body rows 2-18 are a value-snapshot of the refreshed
#    VLOOKUP result (column CB/CC), matching how CM snapshotted the
#    previous pull.
# ---------------------------------------------------------------------
$ws1.Range("CN1").Value = "13-nov"
$ws1.Range("CN1").NumberFormat = "@"

for ($row = 2; $row -le 18; $row++) {
    $src = $ws1.Range("CB$row").Value2
    $dst = $ws1.Range("CN$row")
    $dst.Value = $src
    $dst.NumberFormat = "0"
}

# ---------------------------------------------------------------------
# 3) Leave the cursor where the author left it after this edit.
# ---------------------------------------------------------------------
$null = $ws1.Range("CP8").Select()
